$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextCell 'D2' '79.648.02'
$ws.Range('E2').Value = '  +4.37%  '
Set-TextCell 'D3' '3.204.89'
$ws.Range('E3').Value = '  +5.40%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextCell 'D5' '211.44'
$ws.Range('E5').Value = '  +6.65%  '
Set-TextCell 'D6' '640.27'
$ws.Range('E6').Value = '  +3.64%  '
Set-TextCell 'D7' '1.00'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +19.28%  '
Set-TextCell 'D9' '0.603'
$ws.Range('E9').Value = '  +9.96%  '
Set-TextCell 'D10' '3.207.24'
$ws.Range('E10').Value = '  +5.45%  '
Set-TextCell 'D11' '0.599'
$ws.Range('E11').Value = '  +37.12%  '
Set-TextCell 'D12' '0.0000258'
$ws.Range('E12').Value = '  +33.94%  '
$ws.Range('E13').Value = '  +3.36%  '
Set-TextCell 'D14' '5.43'
$ws.Range('E14').Value = '  +2.89%  '
Set-TextCell 'D15' '3.800.33'
$ws.Range('E15').Value = '  +5.67%  '
Set-TextCell 'D16' '32.19'
$ws.Range('E16').Value = '  +11.67%  '
Set-TextCell 'D17' '79.498.57'
$ws.Range('E17').Value = '  +4.24%  '
Set-TextCell 'D18' '3.204.40'
$ws.Range('E18').Value = '  +5.21%  '
Set-TextCell 'D19' '14.64'
$ws.Range('E19').Value = '  +8.34%  '
$ws.Range('E20').Value = '  +4.83%  '
$ws.Range('E21').Value = '  +26.07%  '
Set-TextCell 'D22' '438.99'
$ws.Range('E22').Value = '  +14.65%  '
Set-TextCell 'D23' '5.28'
$ws.Range('E23').Value = '  +20.92%  '
Set-TextCell 'D24' '4.84'
$ws.Range('E24').Value = '  +11.75%  '
Set-TextCell 'D25' '3.372.95'
$ws.Range('E25').Value = '  +6.19%  '
Set-TextCell 'D26' '77.47'
$ws.Range('E26').Value = '  +6.99%  '
Set-TextCell 'D27' '10.87'
$ws.Range('E27').Value = '  +11.40%  '
Set-TextCell 'D28' '1.00'
$ws.Range('E28').Value = '  +0.08%  '
Set-TextCell 'D29' '0.0000124'
$ws.Range('E29').Value = '  +15.07%  '
Set-TextCell 'D30' '9.17'
$ws.Range('E30').Value = '  +11.09%  '
Set-TextCell 'D31' '0.999'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D32' '563.22'
$ws.Range('E32').Value = '  +14.09%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell 'D33' '1.53'
$ws.Range('E33').Value = '  +10.19%  '
Set-TextCell 'D34' '0.158'
$ws.Range('E34').Value = '  +31.95%  '
Set-TextCell 'D35' '2.05'
$ws.Range('E35').Value = '  +7.00%  '
Set-TextCell 'D36' '23.11'
$ws.Range('E36').Value = '  +12.45%  '
$ws.Range('E37').Value = '  +18.30%  '
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('E39').Value = '  +9.68%  '
Set-TextCell 'D40' '163.58'
$ws.Range('E40').Value = '  +0.80%  '
Set-TextCell 'D42' '5.70'
$ws.Range('E42').Value = '  +11.54%  '
Set-TextCell 'D43' '192.74'
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('E44').Value = '  +0.03%  '
Set-TextCell 'D45' '1.83'
$ws.Range('E45').Value = '  +11.95%  '
Set-TextCell 'D46' '0.802'
$ws.Range('E46').Value = '  +2.03%  '
$ws.Range('B47').Value = 'ImmutableX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D47' '1.34'
$ws.Range('E47').Value = '  +7.66%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell 'D48' '2.68'
$ws.Range('E48').Value = '  +10.56%  '
Set-TextCell 'D49' '43.37'
$ws.Range('E49').Value = '  +4.03%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 'D50' '25.85'
$ws.Range('E50').Value = '  +16.48%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D51' '0.642'
$ws.Range('E51').Value = '  +7.28%  '
